# Update "想去人数" (interested-attendee count) column F on the two sheets
# that list these exhibitions: "展览" (sheet1) and "全部类型" (sheet4).
# "演出" and "本地生活" are untouched by this edit.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6067
$ws1.Range("F5").Value = 365
$ws1.Range("F9").Value = 45
$ws1.Range("F10").Value = 63
$ws1.Range("F12").Value = 143
$ws1.Range("F13").Value = 348
$ws1.Range("F14").Value = 451
$ws1.Range("F15").Value = 3066
$ws1.Range("F17").Value = 175
$ws1.Range("F18").Value = 1716
$ws1.Range("F19").Value = 17

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6067
$ws4.Range("F5").Value = 365
$ws4.Range("F10").Value = 45
$ws4.Range("F11").Value = 63
$ws4.Range("F13").Value = 143
$ws4.Range("F14").Value = 348
$ws4.Range("F15").Value = 451
$ws4.Range("F16").Value = 3066
$ws4.Range("F18").Value = 175
$ws4.Range("F19").Value = 1716
$ws4.Range("F20").Value = 17
